$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: the expired cert has aged one more day since the last run
$ws.Range("E5").Value = "❌ EXPIRED 3777 days ago"

# Row 9: cert expiry date moved out a year (no longer expiring soon),
# so the warning status is cleared
$ws.Range("B9").Style = "Normal"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "2026-09-04"
$ws.Range("B9").Style = "Normal"
$ws.Range("E9").Value = ""

# Row 10: cert expiry date changed (no longer expiring soon),
# so the warning status is cleared
$ws.Range("B10").Style = "Normal"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "2025-11-11"
$ws.Range("B10").Style = "Normal"
$ws.Range("E10").Value = ""
